# Add a new weekly ranking sheet: magapoke_2026-01-28
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "magapoke_2026-01-28"

$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"

$titles = @(
    'ブルーロック',
    'WIND BREAKER',
    'ギルティサークル',
    '東京卍リベンジャーズ',
    'ベイビーステップ',
    '君が僕らを悪魔と呼んだ頃',
    '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！',
    '島耕作',
    'イレギュラーズ',
    'ガチアクタ',
    '黒猫と魔女の教室',
    '南海トラフ巨大地震',
    '魔女と傭兵',
    'となりの黒川さん',
    '十字架のろくにん',
    '魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～',
    '転生したら第七王子だったので、気ままに魔術を極めます',
    'ハードワーカー中田',
    '異世界ウォーキング',
    '黄昏町プリズナーズ',
    'ひゃくえむ。',
    'K-9~警視庁公安部公安第9課異能対策係~',
    '【爆アド】生まれた直後から最強悪霊と脳内バトルしてたら魔力量が測定可能域を超えてました〜悪憑の子の謙虚な覇道〜',
    'ナキナギ',
    '転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～',
    'ドラハチ',
    'アルキメデスの大戦',
    'デッドアカウント',
    '蒼く染めろ',
    '幼馴染とはラブコメにならない',
    '愛妻の裏アカ',
    '味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す',
    '追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜',
    'せいぶつ部の田辺くん',
    'おやすみ ふみさん',
    'FAIRY TAIL 100 YEARS QUEST',
    'さわらないで小手指くん',
    'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
    'ハナバス　苔石花江のバスケ論',
    '普通の本はありません！',
    '辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜',
    'いじめるヤバイ奴',
    'グラぱらっ！',
    '春くらり',
    '食糧人類-Starving Anonymous-',
    '降り積もれ孤独な死よ',
    '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
    'Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。',
    'それがメイドのカンナです',
    '田んぼで拾った女騎士、田舎で俺の嫁だと思われている',
    '屋根の下のアルテミス',
    '不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～',
    '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
    '生徒と恋はできません！',
    'ジュミドロ',
    'インフェクション',
    '東京卍リベンジャーズ～場地圭介からの手紙～',
    'ストーカー行為がバレて人生終了男',
    'MYS',
    '金田一少年の事件簿外伝 犯人たちの事件簿',
    'DAYS外伝',
    '日本語が話せないロシア人美少女転入生が頼れるのは、多言語マスターの俺1人',
    'デスティニーラバーズ',
    '時々ボソッとロシア語でデレる隣のアーリャさん',
    '阿武ノーマル',
    '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
    '白鳥運子は31画',
    '魁の花巫女',
    '念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！',
    '劣等人の魔剣使い　スキルボードを駆使して最強に至る',
    'ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～',
    'アオバノバスケ',
    'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜',
    'リスナーに騙されてダンジョンの最下層から脱出RTAすることになった',
    '不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～',
    'ぼくのアデリア',
    '鳴るさんだぁ',
    'ともだちづくり',
    'シャングリラ・フロンティア～クソゲーハンター、神ゲーに挑まんとす～',
    '五輪の女神さま 〜なでしこ寮のメダルごはん〜',
    '復讐の教科書',
    '我間乱 ―修羅―',
    '追放されなかった男　～二度目の人生は土下座から始まりました～',
    '勇者と呼ばれた後に　―そして無双男は家族を創る―',
    'おくることば',
    'Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜',
    'なれの果ての僕ら',
    'ヒプノシスマイク -Division Rap Battle- side D.H & B.A.T',
    '母という呪縛 娘という牢獄',
    'ヒロインは絶望しました。',
    '最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～',
    '殺人猟団 ‐マッドメン‐',
    '可愛いだけじゃない式守さん',
    '人間消失',
    '触手魔術師の成り上がり',
    '明智ナンバーワン',
    '中華一番！極',
    '絶望集落',
    '東京ネオンスキャンダル',
    'はじめの一歩',
)

for ($i = 0; $i -lt $titles.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $titles[$i]
}
